$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction rows (rows 2-15), replacing the old rows 2-9 and extending the table.
$data = @(
    @(118, "AGUA CRISTAL BOT X 300 ML", 1, 825, "efectivo", "2025-06-21 00:06:24"),
    @(119, "Desodorante Rexona Men V8 Roll On X 30ml", 1, 3000, "efectivo", "2025-06-21 00:06:24"),
    @(120, "CEPILLO DENTAL COLGATE ULTRA PREMIER", 1, 2050, "efectivo", "2025-06-21 00:06:24"),
    @(121, "tomates", 1, 5000, "efectivo", "2025-06-21 00:06:24"),
    @(122, "Cocacola", 1, 5000, "efectivo", "2025-06-21 01:39:44"),
    @(123, "Desodorante para Pies Rexona Efficient Original 55 G", 1, 8300, "efectivo", "2025-06-21 01:39:44"),
    @(124, "FIBER PRO CLEAN 450g", 2, 30000, "efectivo", "2025-06-21 01:39:44"),
    @(125, "Cocacola", 1, 5000, "efectivo", "2025-06-21 01:58:23"),
    @(126, "AGUA CRISTAL BOT X 300 ML", 1, 825, "efectivo", "2025-06-21 01:58:23"),
    @(127, "Manzana", 1, 2000, "efectivo", "2025-06-21 01:58:23"),
    @(128, "Cocacola", 1, 5000, "efectivo", "2025-06-21 01:59:20"),
    @(129, "Cocacola", 1, 5000, "efectivo", "2025-06-21 01:59:47"),
    @(130, "escoba", 1, 1000, "efectivo", "2025-06-21 15:04:55"),
    @(131, "riquillas", 5, 65000, "efectivo", "2025-06-21 15:04:55")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# Row 16 stays blank (gap), summary block now sits at rows 17-20.
$ws.Cells.Item(17, 1).Value = "Resumen del Día"

$ws.Cells.Item(18, 1).Value = "Total Vendido"
$ws.Cells.Item(18, 2).Value = 138000

$ws.Cells.Item(19, 1).Value = "Productos Vendidos"
$ws.Cells.Item(19, 2).Value = 19

$ws.Cells.Item(20, 1).Value = "Número de Ventas"
$ws.Cells.Item(20, 2).Value = 14
